$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new customer data
# (Force A2 to stay text-typed like the original inline string, then
# restore the default "Normal" style so no extra number-format sticks.)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1007311001"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "ЕВРОСТИКС ЕООД"
$ws.Range("C2").Value = 435
$ws.Range("D2").Value = "42.66956517919332,23.38250309228897"
$ws.Range("E2").Value = 42.66956517919332
$ws.Range("F2").Value = 23.38250309228897

# Remove the now-duplicate row 3 entirely
$ws.Rows("3").Delete()
